$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add match results (Skor1/Skor2) for 26.06.2025 matches (rows 12 and 13)
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 7
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1

# Update the active selection on the sheet
$ws.Range("F15").Select()
